# October 2016 data added
# Sheet1 ("Sheffield spending data availability" tracker):
#  - Row 69 (August 2016): Date Field # (C69) now known to be 8.
#  - Row 71 (October 2016): previously marked "Missing", now the author
#    found the data on the council site, so the row is filled in the same
#    way the other "Available" rows are (Available? / Date Field # /
#    Date Format / All in one file / URL), and the "Missing" comment in E71
#    is removed.
#  - A new URL for the October 2016 data is recorded in G71.
#  - Selection/scroll moves along with the newly-edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 69: Date Field # is now known (was blank) ---
$ws.Range("C69").Value = 8

# --- Row 71: October 2016 data was found ---
# Available?
$ws.Range("B71").Style = "Good"
$ws.Range("B71").Value = "Available"

# Date Field #
$ws.Range("C71").Value = 8

# Date Format (wrapped, like the other Date Format cells)
$ws.Range("D71").Style = "Normal"
$ws.Range("D71").WrapText = $true
$ws.Range("D71").Value = "dd/mm/yyyy"

# Comments ("Missing") cell is no longer applicable - remove it entirely
$ws.Range("E71").Clear()

# All in one file
$ws.Range("F71").Value = "… found the data, not sure why I didn't find it before"

# URL
$ws.Range("G71").Value = "https://data.sheffield.gov.uk/Economy/October-2016-Monthly-Payments-To-Suppliers-Over-25/it4i-3itb"

# --- Scroll/selection follows the edited area ---
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$ws.Range("D71").Select()
